$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 2613.125
$ws.Range("I113").Value = 1989.375
$ws.Range("J113").Value = 3236.875
$ws.Range("K113").Value = 1989.375
$ws.Range("L113").Value = 3236.875
$ws.Range("M113").Value = 1264.625
$ws.Range("N113").Value = -9744.875
$ws.Range("H121").Value = 1174.3889
$ws.Range("I121").Value = 897.5
$ws.Range("J121").Value = 1253.5
$ws.Range("K121").Value = 2692.5
$ws.Range("L121").Value = 3760.5
$ws.Range("M121").Value = -945.5
$ws.Range("N121").Value = -7254.5
$ws.Range("H129").Value = 889.15
$ws.Range("I129").Value = 447
$ws.Range("J129").Value = 893.6161499999999
$ws.Range("K129").Value = 1341
$ws.Range("L129").Value = 2680.84845
$ws.Range("M129").Value = 3659
$ws.Range("N129").Value = -12680.84845
$ws.Range("H138").Value = 3657.8901
$ws.Range("I138").Value = 2797.05
$ws.Range("J138").Value = 3900.3804
$ws.Range("K138").Value = 8391.150000000001
$ws.Range("L138").Value = 11701.1412
$ws.Range("M138").Value = -3251.150000000001
$ws.Range("N138").Value = -21981.1412

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2996.3809
$ws.Range("I61").Value = 3148
$ws.Range("J61").Value = 2750
$ws.Range("K61").Value = 3148
$ws.Range("L61").Value = 2750
$ws.Range("M61").Value = -2936
$ws.Range("N61").Value = -3174
$ws.Range("H74").Value = 1851.5
$ws.Range("I74").Value = 1835.6471
$ws.Range("K74").Value = 1835.6471
$ws.Range("M74").Value = -961.6470999999999
$ws.Range("H77").Value = 1851.5
$ws.Range("I77").Value = 1835.6471
$ws.Range("K77").Value = 9178.235499999999
$ws.Range("M77").Value = -4810.235499999999
$ws.Range("H92").Value = 68625
$ws.Range("J92").Value = 68625
$ws.Range("L92").Value = 68625
$ws.Range("N92").Value = -73617
$ws.Range("H95").Value = 29266.4
$ws.Range("J95").Value = 29266.4
$ws.Range("L95").Value = 29266.4
$ws.Range("N95").Value = -34758.4
$ws.Range("H96").Value = 30000
$ws.Range("J96").Value = 30000
$ws.Range("L96").Value = 30000
$ws.Range("N96").Value = -35492
$ws.Range("H132").Value = 7383.075
$ws.Range("I132").Value = 6938.5654
$ws.Range("J132").Value = 7984.4707
$ws.Range("K132").Value = 20815.6962
$ws.Range("L132").Value = 23953.4121
$ws.Range("M132").Value = -18285.6962
$ws.Range("N132").Value = -29013.4121
$ws.Range("H136").Value = 2996.3809
$ws.Range("I136").Value = 3148
$ws.Range("J136").Value = 2750
$ws.Range("K136").Value = 9444
$ws.Range("L136").Value = 8250
$ws.Range("M136").Value = -6894
$ws.Range("N136").Value = -13350

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1770.591
$ws.Range("I134").Value = 1397.0714
$ws.Range("J134").Value = 2424.25
$ws.Range("K134").Value = 4191.2142
$ws.Range("L134").Value = 7272.75
$ws.Range("M134").Value = -1656.2142
$ws.Range("N134").Value = -12342.75

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2132.7896
$ws.Range("I58").Value = 1981.5333
$ws.Range("J58").Value = 2700
$ws.Range("K58").Value = 1981.5333
$ws.Range("L58").Value = 2700
$ws.Range("M58").Value = -1778.5333
$ws.Range("N58").Value = -3106
$ws.Range("H86").Value = 6964.25
$ws.Range("I86").Value = 3928.5
$ws.Range("J86").Value = 10000
$ws.Range("K86").Value = 3928.5
$ws.Range("L86").Value = 10000
$ws.Range("M86").Value = -2805.5
$ws.Range("N86").Value = -12246
$ws.Range("H89").Value = 6964.25
$ws.Range("I89").Value = 3928.5
$ws.Range("J89").Value = 10000
$ws.Range("K89").Value = 19642.5
$ws.Range("L89").Value = 50000
$ws.Range("M89").Value = -14026.5
$ws.Range("N89").Value = -61232
$ws.Range("H94").Value = 1083.3636
$ws.Range("I94").Value = 473.55554
$ws.Range("J94").Value = 1505.5385
$ws.Range("K94").Value = 473.55554
$ws.Range("L94").Value = 1505.5385
$ws.Range("M94").Value = -22.55554000000001
$ws.Range("N94").Value = -2407.5385
$ws.Range("H99").Value = 1731.8235
$ws.Range("I99").Value = 1767
$ws.Range("J99").Value = 1617.5
$ws.Range("K99").Value = 1767
$ws.Range("L99").Value = 1617.5
$ws.Range("M99").Value = -269
$ws.Range("N99").Value = -4613.5
$ws.Range("H126").Value = 1731.8235
$ws.Range("I126").Value = 1767
$ws.Range("J126").Value = 1617.5
$ws.Range("K126").Value = 5301
$ws.Range("L126").Value = 4852.5
$ws.Range("M126").Value = -2831
$ws.Range("N126").Value = -9792.5
$ws.Range("H136").Value = 2132.7896
$ws.Range("I136").Value = 1981.5333
$ws.Range("J136").Value = 2700
$ws.Range("K136").Value = 5944.5999
$ws.Range("L136").Value = 8100
$ws.Range("M136").Value = -3394.5999
$ws.Range("N136").Value = -13200

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 205152.08
$ws.Range("I68").Value = 313203.2
$ws.Range("J68").Value = 1761.7646
$ws.Range("K68").Value = 939609.6000000001
$ws.Range("L68").Value = 5285.293799999999
$ws.Range("M68").Value = -938798.6000000001
$ws.Range("N68").Value = -6907.293799999999
$ws.Range("H71").Value = 205152.08
$ws.Range("I71").Value = 313203.2
$ws.Range("J71").Value = 1761.7646
$ws.Range("K71").Value = 2818828.8
$ws.Range("L71").Value = 15855.8814
$ws.Range("M71").Value = -2814772.8
$ws.Range("N71").Value = -23967.8814
$ws.Range("H97").Value = 1773.5385
$ws.Range("I97").Value = 3000
$ws.Range("J97").Value = 1671.3334
$ws.Range("K97").Value = 9000
$ws.Range("L97").Value = 5014.0002
$ws.Range("M97").Value = -8504
$ws.Range("N97").Value = -6006.0002
$ws.Range("H103").Value = 2352.6667
$ws.Range("I103").Value = 1536.5
$ws.Range("J103").Value = 3005.6
$ws.Range("K103").Value = 4609.5
$ws.Range("L103").Value = 9016.799999999999
$ws.Range("M103").Value = -3730.5
$ws.Range("N103").Value = -10774.8
$ws.Range("H106").Value = 7131.4287
$ws.Range("J106").Value = 7131.4287
$ws.Range("L106").Value = 21394.2861
$ws.Range("N106").Value = -23286.2861
$ws.Range("H131").Value = 20411430
$ws.Range("J131").Value = 21279960
$ws.Range("L131").Value = 63839880
$ws.Range("N131").Value = -63849960

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 121666.664
$ws.Range("J29").Value = 15000
$ws.Range("L29").Value = 15000
$ws.Range("N29").Value = -15580
$ws.Range("H132").Value = 4446.5454
$ws.Range("I132").Value = 3580
$ws.Range("J132").Value = 5168.6665
$ws.Range("K132").Value = 10740
$ws.Range("L132").Value = 15505.9995
$ws.Range("M132").Value = -8210
$ws.Range("N132").Value = -20565.9995

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3040.3
$ws.Range("I7").Value = 2914.7144
$ws.Range("J7").Value = 3333.3333
$ws.Range("K7").Value = 2914.7144
$ws.Range("L7").Value = 3333.3333
$ws.Range("M7").Value = -2802.7144
$ws.Range("N7").Value = -3557.3333
$ws.Range("H122").Value = 21434142
$ws.Range("I122").Value = 35718784
$ws.Range("K122").Value = 107156352
$ws.Range("M122").Value = -107153902
$ws.Range("H126").Value = 3040.3
$ws.Range("I126").Value = 2914.7144
$ws.Range("J126").Value = 3333.3333
$ws.Range("K126").Value = 8744.143199999999
$ws.Range("L126").Value = 9999.999899999999
$ws.Range("M126").Value = -6274.143199999999
$ws.Range("N126").Value = -14939.9999
$ws.Range("H134").Value = 49803.625
$ws.Range("J134").Value = 49803.625
$ws.Range("L134").Value = 49803.625
$ws.Range("N134").Value = -59943.625

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H48").Value = 41666.668
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 41666.668
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 41666.668
$ws.Range("M48").ClearContents()
$ws.Range("N48").Value = -42804.668
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("H136").Value = 2311.125
$ws.Range("I136").Value = 2260.842
$ws.Range("J136").Value = 2384.6155
$ws.Range("K136").Value = 6782.526
$ws.Range("L136").Value = 7153.8465
$ws.Range("M136").Value = -4232.526
$ws.Range("N136").Value = -12253.8465
